# Apply scheduled-runner market data updates to Golem_Profits sheets
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 268.16666
$ws.Range("I19").Value = 359.625
$ws.Range("J19").Value = 85.25
$ws.Range("K19").Value = 359.625
$ws.Range("L19").Value = 85.25
$ws.Range("M19").Value = -184.625
$ws.Range("N19").Value = -435.25
# row 40
$ws.Range("H40").Value = 2733.111
$ws.Range("J40").Value = 3279.8
$ws.Range("L40").Value = 3279.8
$ws.Range("N40").Value = -3629.8
# row 51
$ws.Range("H51").Value = 6067.5454
$ws.Range("I51").Value = 5718.125
$ws.Range("J51").Value = 6999.3335
$ws.Range("K51").Value = 5718.125
$ws.Range("L51").Value = 6999.3335
$ws.Range("M51").Value = -5234.125
$ws.Range("N51").Value = -7967.3335
# row 53
$ws.Range("H53").Value = 344.33334
$ws.Range("I53").Value = 30
$ws.Range("K53").Value = 30
$ws.Range("M53").Value = 607
# row 100
$ws.Range("H100").Value = 2371
$ws.Range("I100").Value = 2134.3
$ws.Range("J100").Value = 3160
$ws.Range("K100").Value = 2134.3
$ws.Range("L100").Value = 3160
$ws.Range("M100").Value = -1593.3
$ws.Range("N100").Value = -4242
# row 132
$ws.Range("H132").Value = 1541.9
$ws.Range("I132").Value = 1379.8889
$ws.Range("K132").Value = 4139.6667
$ws.Range("M132").Value = -1609.6667

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 41
$ws.Range("H41").Value = 588.75
$ws.Range("I41").Value = 588.75
$ws.Range("K41").Value = 588.75
$ws.Range("M41").Value = -174.75
# row 45
$ws.Range("H45").Value = 2077.3333
$ws.Range("J45").Value = 1234
$ws.Range("L45").Value = 1234
$ws.Range("N45").Value = -1988

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 49
$ws.Range("H49").Value = 22500
$ws.Range("J49").Value = 22500
$ws.Range("L49").Value = 22500
$ws.Range("N49").Value = -22978
# row 86
$ws.Range("H86").Value = 3504.5
$ws.Range("I86").Value = 1318.5
$ws.Range("J86").Value = 4597.5
$ws.Range("K86").Value = 1318.5
$ws.Range("L86").Value = 4597.5
$ws.Range("M86").Value = -195.5
$ws.Range("N86").Value = -6843.5
# row 89
$ws.Range("H89").Value = 3504.5
$ws.Range("I89").Value = 1318.5
$ws.Range("J89").Value = 4597.5
$ws.Range("K89").Value = 6592.5
$ws.Range("L89").Value = 22987.5
$ws.Range("M89").Value = -976.5
$ws.Range("N89").Value = -34219.5
# row 99
$ws.Range("H99").Value = 3486.6155
$ws.Range("I99").Value = 3486.6155
$ws.Range("K99").Value = 3486.6155
$ws.Range("M99").Value = -1988.6155
# row 105
$ws.Range("H105").Value = 922.5454999999999
$ws.Range("I105").Value = 813.7
$ws.Range("K105").Value = 813.7
$ws.Range("M105").Value = 933.3

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 1596.25
$ws.Range("I16").Value = 1461.6666
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1461.6666
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1174.6666
$ws.Range("N16").Value = -2574
# row 57
$ws.Range("H57").Value = 48000
$ws.Range("J57").Value = 48000
$ws.Range("L57").Value = 48000
$ws.Range("N57").Value = -49120
# row 58
$ws.Range("H58").Value = 1500
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 1000
$ws.Range("M58").Value = -797
# row 69
$ws.Range("H69").Value = 35639.4
$ws.Range("J69").Value = 43999.25
$ws.Range("L69").Value = 43999.25
$ws.Range("N69").Value = -45497.25
# row 72
$ws.Range("H72").Value = 35639.4
$ws.Range("J72").Value = 43999.25
$ws.Range("L72").Value = 131997.75
$ws.Range("N72").Value = -139485.75
# row 105
$ws.Range("H105").Value = 445.75
$ws.Range("I105").Value = 445.75
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 445.75
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1301.25
$ws.Range("N105").ClearContents()
# row 107
$ws.Range("H107").Value = 677.4
$ws.Range("I107").Value = 647
$ws.Range("K107").Value = 647
$ws.Range("M107").Value = 1273
# row 113
$ws.Range("H113").Value = 1596.25
$ws.Range("I113").Value = 1461.6666
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1461.6666
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 708.3334
$ws.Range("N113").Value = -6340
# row 132
$ws.Range("H132").Value = 1542.3334
$ws.Range("I132").Value = 1513.5
$ws.Range("K132").Value = 4540.5
$ws.Range("M132").Value = -2010.5
# row 136
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 12
$ws.Range("H12").Value = 217.66667
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 326
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 978
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -1324
# row 23
$ws.Range("H23").Value = 126.92857
$ws.Range("J23").Value = 161.3
$ws.Range("L23").Value = 483.9
$ws.Range("N23").Value = -953.9000000000001
# row 34
$ws.Range("H34").Value = 3329.1667
$ws.Range("J34").Value = 3382.9575
$ws.Range("L34").Value = 10148.8725
$ws.Range("N34").Value = -10316.8725
# row 39
$ws.Range("H39").Value = 2333.3333
$ws.Range("J39").Value = 2333.3333
$ws.Range("L39").Value = 6999.999899999999
$ws.Range("N39").Value = -7587.999899999999
# row 55
$ws.Range("H55").Value = 3534.984
$ws.Range("J55").Value = 3735.2678
$ws.Range("L55").Value = 11205.8034
$ws.Range("N55").Value = -11559.8034
# row 80
$ws.Range("H80").Value = 4982.6665
$ws.Range("I80").Value = 4983
$ws.Range("K80").Value = 14949
$ws.Range("M80").Value = -14013
# row 83
$ws.Range("H83").Value = 4982.6665
$ws.Range("I83").Value = 4983
$ws.Range("K83").Value = 44847
$ws.Range("M83").Value = -40167
# row 107
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
# row 139
$ws.Range("H139").Value = 57325.445

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 62
$ws.Range("H62").Value = 29999
$ws.Range("J62").Value = 29999
$ws.Range("L62").Value = 29999
$ws.Range("N62").Value = -31371
# row 65
$ws.Range("H65").Value = 29999
$ws.Range("J65").Value = 29999
$ws.Range("L65").Value = 89997
$ws.Range("N65").Value = -96861
# row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 1772.8334
$ws.Range("I22").Value = 1344
$ws.Range("J22").Value = 1987.25
$ws.Range("K22").Value = 1344
$ws.Range("L22").Value = 1987.25
$ws.Range("M22").Value = -1049
$ws.Range("N22").Value = -2577.25
# row 27
$ws.Range("H27").Value = 1772.8334
$ws.Range("I27").Value = 1344
$ws.Range("J27").Value = 1987.25
$ws.Range("K27").Value = 1344
$ws.Range("L27").Value = 1987.25
$ws.Range("M27").Value = -1237
$ws.Range("N27").Value = -2201.25
# row 82
$ws.Range("H82").Value = 3163
$ws.Range("I82").Value = 2136.7144
$ws.Range("K82").Value = 2136.7144
$ws.Range("M82").Value = -1775.7144
# row 85
$ws.Range("H85").Value = 3163
$ws.Range("I85").Value = 2136.7144
$ws.Range("K85").Value = 2136.7144
$ws.Range("M85").Value = -888.7143999999998
# row 103
$ws.Range("H103").Value = 5999.75
$ws.Range("J103").Value = 5999.75
$ws.Range("L103").Value = 5999.75
$ws.Range("N103").Value = -8343.75

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 136
$ws.Range("H136").Value = 1404.3334
$ws.Range("I136").Value = 1404.3334
$ws.Range("K136").Value = 4213.0002
$ws.Range("M136").Value = -1663.0002
